$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find last used row based on column A (Beteckning) to know data extent.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C ("Förändrad") holds a date serial that must move from 45178 to 45179
# for every data row (row 2 through the last used row).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
